$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update player "DIEGO" (row 6) stats: pontuacao, moedas, partidas, mediapontos
$ws.Range("C6").Value = 100
$ws.Range("D6").Value = 100
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 50
